# Natmi following Dr Hou advice
# Update C1qb-Lrp1 LR-pair sheet: refresh stats for ECs and add sCs as a new sending cluster
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove previous data rows (rows 2-4) so the sheet can be rebuilt with the new data (rows 2-7)
$ws.Range("A2:T4").ClearContents()

# Row 2: ECs -> ECs (C1qb-Lrp1)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "C1qb"
$ws.Range("C2").Value = "Lrp1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 11.81087266666667
$ws.Range("H2").Value = 35.432618
$ws.Range("I2").Value = 0.9915571911324677
$ws.Range("J2").Value = 0.9915571911324677
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 21.08181366666667
$ws.Range("N2").Value = 63.245441
$ws.Range("O2").Value = 0.0571606014598545
$ws.Range("P2").Value = 0.0571606014598545
$ws.Range("Q2").Value = 248.9946167993931
$ws.Range("R2").Value = 2240.951551194538
$ws.Range("S2").Value = 0.05667800542697576
$ws.Range("T2").Value = 0.05667800542697576

# Row 3: ECs -> FAPs (C1qb-Lrp1)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "C1qb"
$ws.Range("C3").Value = "Lrp1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 11.81087266666667
$ws.Range("H3").Value = 35.432618
$ws.Range("I3").Value = 0.9915571911324677
$ws.Range("J3").Value = 0.9915571911324677
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 301.6001486666667
$ws.Range("N3").Value = 904.800446
$ws.Range("O3").Value = 0.8177496571571792
$ws.Range("P3").Value = 0.8177496571571792
$ws.Range("Q3").Value = 3562.160952149736
$ws.Range("R3").Value = 32059.44856934762
$ws.Range("S3").Value = 0.810845553100311
$ws.Range("T3").Value = 0.810845553100311

# Row 4: ECs -> sCs (C1qb-Lrp1)
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "C1qb"
$ws.Range("C4").Value = "Lrp1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 11.81087266666667
$ws.Range("H4").Value = 35.432618
$ws.Range("I4").Value = 0.9915571911324677
$ws.Range("J4").Value = 0.9915571911324677
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 46.13524966666667
$ws.Range("N4").Value = 138.405749
$ws.Range("O4").Value = 0.1250897413829664
$ws.Range("P4").Value = 0.1250897413829664
$ws.Range("Q4").Value = 544.8975592578759
$ws.Range("R4").Value = 4904.078033320882
$ws.Range("S4").Value = 0.1240336326051809
$ws.Range("T4").Value = 0.1240336326051809

# Row 5: sCs -> ECs (C1qb-Lrp1)
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "C1qb"
$ws.Range("C5").Value = "Lrp1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.100566
$ws.Range("H5").Value = 0.301698
$ws.Range("I5").Value = 0.008442808867532263
$ws.Range("J5").Value = 0.008442808867532263
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 21.08181366666667
$ws.Range("N5").Value = 63.245441
$ws.Range("O5").Value = 0.0571606014598545
$ws.Range("P5").Value = 0.0571606014598545
$ws.Range("Q5").Value = 2.120113673202
$ws.Range("R5").Value = 19.081023058818
$ws.Range("S5").Value = 0.0004825960328787372
$ws.Range("T5").Value = 0.0004825960328787372

# Row 6: sCs -> FAPs (C1qb-Lrp1)
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "C1qb"
$ws.Range("C6").Value = "Lrp1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.100566
$ws.Range("H6").Value = 0.301698
$ws.Range("I6").Value = 0.008442808867532263
$ws.Range("J6").Value = 0.008442808867532263
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 301.6001486666667
$ws.Range("N6").Value = 904.800446
$ws.Range("O6").Value = 0.8177496571571792
$ws.Range("P6").Value = 0.8177496571571792
$ws.Range("Q6").Value = 30.330720550812
$ws.Range("R6").Value = 272.976484957308
$ws.Range("S6").Value = 0.0069041040568681
$ws.Range("T6").Value = 0.0069041040568681

# Row 7: sCs -> sCs (C1qb-Lrp1)
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "C1qb"
$ws.Range("C7").Value = "Lrp1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.100566
$ws.Range("H7").Value = 0.301698
$ws.Range("I7").Value = 0.008442808867532263
$ws.Range("J7").Value = 0.008442808867532263
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 46.13524966666667
$ws.Range("N7").Value = 138.405749
$ws.Range("O7").Value = 0.1250897413829664
$ws.Range("P7").Value = 0.1250897413829664
$ws.Range("Q7").Value = 4.639637517978001
$ws.Range("R7").Value = 41.75673766180201
$ws.Range("S7").Value = 0.001056108777785426
$ws.Range("T7").Value = 0.001056108777785426
